$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "Id" column (A2:A11) from numeric 0-9 to uint-style unique string
# identifiers S001-S010, per the "uint unique ID support" change.
$ws.Range("A2").Value = "S001"
$ws.Range("A3").Value = "S002"
$ws.Range("A4").Value = "S003"
$ws.Range("A5").Value = "S004"
$ws.Range("A6").Value = "S005"
$ws.Range("A7").Value = "S006"
$ws.Range("A8").Value = "S007"
$ws.Range("A9").Value = "S008"
$ws.Range("A10").Value = "S009"
$ws.Range("A11").Value = "S010"
